# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" sheet (right after "总计") with the latest quarterly
# fund-holdings snapshot, and records the new quarter in the "总计" summary
# sheet (new top row, existing rows shift down by one).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert the 2022-Q4 totals as the new first data
#    row, pushing the previously-existing rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Give the brand-new row (row 7) the same formatting as column A already has
# on the other data rows before we start overwriting values.
$total.Cells.Item(6, 1).Copy()
$total.Cells.Item(7, 1).PasteSpecial(-4122)

$totalRows = @(
    @(0, "2022-Q4", 4,  0.6899999999999999),
    @(1, "2022-Q2", 7,  0.5600000000000001),
    @(2, "2021-Q4", 7,  0.4),
    @(3, "2021-Q2", 2,  2.14),
    @(4, "2021-Q1", 13, 4.29),
    @(5, "2020-Q4", 4,  1.95)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2. New "2022-Q4" worksheet, placed right after "总计", carrying the same
#    layout/styling as the other quarterly sheets (e.g. "2022-Q2").
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

# Reuse the header-row and column-A formatting from an existing quarter sheet.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking data stored as TEXT (matching
# the source workbook's inlineStr cells) - a leading apostrophe forces text
# the way typing it directly into Excel would, without touching NumberFormat.
$q4Rows = @(
    @(0, "001487", "宝盈优势产业灵活配置混合A", "10.62", "94.48", "2.72", "0.2889", 10),
    @(1, "013895", "宝盈成长精选混合A",         "8.59",  "94.68", "2.71", "0.2328", 10),
    @(2, "012771", "宝盈优势产业灵活配置混合C", "3.19",  "94.48", "2.72", "0.0868", 10),
    @(3, "013896", "宝盈成长精选混合C",         "3.06",  "94.68", "2.71", "0.0829", 10)
)

for ($i = 0; $i -lt $q4Rows.Count; $i++) {
    $r = $i + 2
    $row = $q4Rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}
